# Upload new version with timestamp
#
# The workbook's "عدد التعاملات" (transactions) column stores a
# "days:hours" style counter. This commit refreshes that counter:
#   - rows that used to read "0:0" now read "-1:0"
#   - the QUICK NAIL LOTION row (row 9) gets its own distinct value "-23:0"
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7 and 8 (ORGASOL LIGHT CREAM, PRISBRINA CAPS) shared the "0:0" text;
# it is refreshed to "-1:0".
$ws.Range("H7").Value = "-1:0"
$ws.Range("H8").Value = "-1:0"

# Row 9 (QUICK NAIL LOTION) gets its own new value "-23:0".
$ws.Range("H9").Value = "-23:0"
